$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number by Excel;
# force text storage (matching the source inlineStr "t" type) then restore the default style.
$ws.Range('D2').Value = '70.076.55'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '3.619.68'
$ws.Range('E3').Value = '  +3.49%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '195.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.626'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.207'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.10%  '
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.96'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('E12').Value = '  +1.38%  '
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').Value = '4.190.14'
$ws.Range('E14').Value = '  +3.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '590.85'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.23'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').Value = '70.237.77'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').Value = '3.623.95'
$ws.Range('E19').Value = '  +3.44%  '
$ws.Range('E20').Value = '  +1.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.994'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.72'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.80%  '
$ws.Range('E23').Value = '  +2.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '102.46'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.93%  '
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('E26').Value = '  -1.60%  '
$ws.Range('E27').Value = '  -1.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.61'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.86'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.46'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.15'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('E32').Value = '  -2.84%  '
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.24'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('E35').Value = '  +12.07%  '
$ws.Range('D36').Value = '3.957.17'
$ws.Range('E36').Value = '  +5.84%  '
$ws.Range('E37').Value = '  +5.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '528.45'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.60%  '
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.22'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.62%  '
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('E42').Value = '  +0.49%  '
$ws.Range('E43').Value = '  -1.90%  '
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.28%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.141'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.80%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.34'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.90%  '
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('E50').Value = '  +5.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.32'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.29%  '
